$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (same across every model row, columns B..Q)
$values = @(
    0.9994384710386303,
    0.9988584094093185,
    0.9999800373125493,
    0.9998395362326661,
    0.9999729537302476,
    0.0005241626371659963,
    0.001065624706365756,
    0.00001882471861840005,
    0.00003589521605598040,
    0.00002735996733719023,
    0.0009694867605113856,
    0.02289459842770771,
    1.002695339014575,
    0.02386926888112749,
    73.107417091058,
    108.4548160122358
)

# Columns B (2) through Q (17)
$firstCol = 2
$lastRow = 26

for ($row = 2; $row -le $lastRow; $row++) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $firstCol + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
